$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) The trailing columns U:AD (rows 1-19) held a duplicated/unused block of
#    header numbers + pair labels. Remove them so the sheet shrinks back
#    down to A1:T19 (before the new rows are appended below).
$ws.Range("U1:AD19").Delete()

# 2) Row 2 is the HKL header row (C2:T2). The underlying label ordering was
#    reshuffled; update the displayed text to match the new order.
$ws.Range("C2").Value = "[2, 2, 0]"
$ws.Range("D2").Value = "[2, 0, 0]"
$ws.Range("E2").Value = "[4, 0, 0]"
$ws.Range("F2").Value = "[2, 1, 1]"
$ws.Range("G2").Value = "[3, 2, 1]"
$ws.Range("H2").Value = "[2, 2, 2]"
$ws.Range("I2").Value = "[3, 1, 0]"
$ws.Range("J2").Value = "[1, 1, 0]"
$ws.Range("K2").Value = "1Pair-A"
$ws.Range("L2").Value = "1Pair-B"
$ws.Range("M2").Value = "2Pairs-A"
$ws.Range("N2").Value = "2Pairs-B"
$ws.Range("O2").Value = "3Pairs-A"
$ws.Range("P2").Value = "3Pairs-B"
$ws.Range("Q2").Value = "3Pairs-C"
$ws.Range("R2").Value = "4Pairs"
$ws.Range("S2").Value = "5A4F"
$ws.Range("T2").Value = "MaxUnique"

# 3) Rows 16-19 used to be the "HexGrid-90degTilt*degRes" simulations. They
#    are renamed to the new "Holden*" scheme (same HKL index / row data).
$ws.Range("B16").Value = "Holden2.5"
$ws.Range("B17").Value = "Holden5"
$ws.Range("B18").Value = "Holden10"
$ws.Range("B19").Value = "Holden15"

# 4) The displaced HexGrid rows are appended as new rows 20-23, keeping the
#    same column layout (HKL index in A, label in B, 1's across C:T).
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "HexGrid-90degTilt2.5degRes"
$ws.Range("C20:T20").Value = 1

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C21:T21").Value = 1

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "HexGrid-90degTilt10degRes"
$ws.Range("C22:T22").Value = 1

$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "HexGrid-90degTilt15degRes"
$ws.Range("C23:T23").Value = 1

# Match the bold + bordered + centered/top-aligned style used by column A
# elsewhere in the sheet (same as A3:A19) for the new rows' A cells.
$newA = $ws.Range("A20:A23")
$newA.Font.Bold = $true
$newA.Borders.LineStyle = 1
$newA.HorizontalAlignment = -4108
$newA.VerticalAlignment = -4160
